$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.673.22"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").Value = "2.223.73"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'274.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.29%  "

$ws.Range("D6").Value = "'86.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.69%  "

$ws.Range("D7").Value = "'0.614"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.10%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -0.59%  "

$ws.Range("D10").Value = "'45.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.47%  "

$ws.Range("D11").Value = "'0.0917"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").Value = "'7.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.39%  "

$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("D14").Value = "2.556.38"
$ws.Range("E14").Value = "  -1.08%  "

$ws.Range("D15").Value = "'14.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").Value = "2.205.09"
$ws.Range("E16").Value = "  -1.33%  "

$ws.Range("D17").Value = "'0.786"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.29%  "

$ws.Range("D18").Value = "43.593.20"
$ws.Range("E18").Value = "  -1.08%  "

$ws.Range("E19").Value = "  -1.58%  "

$ws.Range("D20").Value = "'70.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("E21").Value = "  -2.17%  "

$ws.Range("D22").Value = "'2.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").Value = "'232.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.42%  "

$ws.Range("D24").Value = "'8.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.47%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "'2.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.19%  "

$ws.Range("D27").Value = "'10.81"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").Value = "'3.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.45%  "

$ws.Range("D29").Value = "'2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.60%  "

$ws.Range("D30").Value = "'39.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.78%  "

$ws.Range("D31").Value = "'172.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("D32").Value = "'0.0907"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.86%  "

$ws.Range("D33").Value = "'20.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").Value = "'5.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  -5.20%  "

$ws.Range("D37").Value = "'0.0353"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.66%  "

$ws.Range("D38").Value = "'4.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.10%  "

$ws.Range("D39").Value = "'3.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.04%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'12.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.12%  "

$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'2.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.83%  "

$ws.Range("D42").Value = "'63.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "'5.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.12%  "

$ws.Range("D45").Value = "'8.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.87%  "

$ws.Range("D46").Value = "'0.0984"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").Value = "'100.04"
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = "  +2.81%  "

$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("D50").Value = "'0.426"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.80%  "

$ws.Range("E51").Value = "  -2.81%  "
